$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 195 (shifts existing rows 195:206 down to 196:207)
$ws.Rows.Item(195).Insert()

# Populate the newly inserted row 195 with the new data record
$row = 195
$ws.Cells.Item($row, 1).Value  = 5
$ws.Cells.Item($row, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item($row, 3).Value  = "Maule"
$ws.Cells.Item($row, 4).Value  = 45166
$ws.Cells.Item($row, 5).Value  = 7
$ws.Cells.Item($row, 6).Value  = 100112001
$ws.Cells.Item($row, 7).Value  = "Berenjena"
$ws.Cells.Item($row, 8).Value  = "Sin especificar"
$ws.Cells.Item($row, 9).Value  = "Primera"
$ws.Cells.Item($row, 10).Value = 200
$ws.Cells.Item($row, 11).Value = 10000
$ws.Cells.Item($row, 12).Value = 10000
$ws.Cells.Item($row, 13).Value = 10000
$ws.Cells.Item($row, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value = 200
$ws.Cells.Item($row, 17).Value = 50
$ws.Cells.Item($row, 18).Value = "Hortaliza"
